$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 4586.75
$ws.Range("J17").Value = 4586.75
$ws.Range("L17").Value = 13760.25
$ws.Range("N17").Value = -14096.25
$ws.Range("H41").Value = 213.21053
$ws.Range("I41").Value = 125.53333
$ws.Range("J41").Value = 542
$ws.Range("K41").Value = 125.53333
$ws.Range("L41").Value = 542
$ws.Range("M41").Value = 314.46667
$ws.Range("N41").Value = -1422
$ws.Range("H80").Value = 8379.52
$ws.Range("I80").Value = 6891.1333
$ws.Range("J80").Value = 10612.1
$ws.Range("K80").Value = 20673.3999
$ws.Range("L80").Value = 31836.3
$ws.Range("M80").Value = -19675.3999
$ws.Range("N80").Value = -33832.3
$ws.Range("H81").Value = 0
$ws.Range("J81").Value = 0
$ws.Range("L81").Value = 0
$ws.Range("N81").ClearContents()
$ws.Range("H83").Value = 8379.52
$ws.Range("I83").Value = 6891.1333
$ws.Range("J83").Value = 10612.1
$ws.Range("K83").Value = 62020.1997
$ws.Range("L83").Value = 95508.90000000001
$ws.Range("M83").Value = -57028.1997
$ws.Range("N83").Value = -105492.9
$ws.Range("H84").Value = 0
$ws.Range("J84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("N84").ClearContents()
$ws.Range("H92").Value = 1208.4375
$ws.Range("I92").Value = 865.8461
$ws.Range("K92").Value = 865.8461
$ws.Range("M92").Value = 382.1539
$ws.Range("H98").Value = 1371.875
$ws.Range("I98").Value = 1416.3043
$ws.Range("J98").Value = 350
$ws.Range("K98").Value = 1416.3043
$ws.Range("L98").Value = 350
$ws.Range("M98").Value = 81.69569999999999
$ws.Range("N98").Value = -3346
$ws.Range("H122").Value = 1371.875
$ws.Range("I122").Value = 1416.3043
$ws.Range("J122").Value = 350
$ws.Range("K122").Value = 4248.9129
$ws.Range("L122").Value = 1050
$ws.Range("M122").Value = -1798.9129
$ws.Range("N122").Value = -5950
$ws.Range("H138").Value = 6889.143
$ws.Range("I138").Value = 6942.6
$ws.Range("J138").Value = 6817.8667
$ws.Range("K138").Value = 20827.8
$ws.Range("L138").Value = 20453.6001
$ws.Range("M138").Value = -15687.8
$ws.Range("N138").Value = -30733.6001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 10645.107
$ws.Range("I32").Value = 9135.367
$ws.Range("J32").Value = 34499
$ws.Range("K32").Value = 9135.367
$ws.Range("L32").Value = 34499
$ws.Range("M32").Value = -8848.367
$ws.Range("N32").Value = -35073
$ws.Range("H45").Value = 86538.586
$ws.Range("I45").Value = 119813.414
$ws.Range("K45").Value = 119813.414
$ws.Range("M45").Value = -119436.414
$ws.Range("H61").Value = 1889.3214
$ws.Range("I61").Value = 1889.3214
$ws.Range("K61").Value = 1889.3214
$ws.Range("M61").Value = -1677.3214
$ws.Range("H97").Value = 884.129
$ws.Range("I97").Value = 939.1923
$ws.Range("K97").Value = 939.1923
$ws.Range("M97").Value = -443.1923
$ws.Range("H123").Value = 20166.666
$ws.Range("J123").Value = 20166.666
$ws.Range("L123").Value = 20166.666
$ws.Range("N123").Value = -29966.666
$ws.Range("H132").Value = 2766.484
$ws.Range("I132").Value = 2800.0356
$ws.Range("J132").Value = 2453.3333
$ws.Range("K132").Value = 8400.106800000001
$ws.Range("L132").Value = 7359.999899999999
$ws.Range("M132").Value = -5870.106800000001
$ws.Range("N132").Value = -12419.9999
$ws.Range("H136").Value = 1889.3214
$ws.Range("I136").Value = 1889.3214
$ws.Range("K136").Value = 5667.9642
$ws.Range("M136").Value = -3117.9642

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H2").Value = 0
$ws.Range("J2").Value = 0
$ws.Range("L2").Value = 0
$ws.Range("N2").ClearContents()
$ws.Range("H86").Value = 1954.7273
$ws.Range("I86").Value = 1844
$ws.Range("K86").Value = 1844
$ws.Range("M86").Value = -721
$ws.Range("H89").Value = 1954.7273
$ws.Range("I89").Value = 1844
$ws.Range("K89").Value = 9220
$ws.Range("M89").Value = -3604
$ws.Range("H134").Value = 3916.1365
$ws.Range("I134").Value = 3323.2104
$ws.Range("K134").Value = 9969.6312
$ws.Range("M134").Value = -7434.6312

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 8498.637000000001
$ws.Range("I58").Value = 5331.3335
$ws.Range("J58").Value = 12299.4
$ws.Range("K58").Value = 5331.3335
$ws.Range("L58").Value = 12299.4
$ws.Range("M58").Value = -5128.3335
$ws.Range("N58").Value = -12705.4
$ws.Range("H64").Value = 86666.664
$ws.Range("J64").Value = 86666.664
$ws.Range("L64").Value = 86666.664
$ws.Range("N64").Value = -87162.664
$ws.Range("H67").Value = 86666.664
$ws.Range("J67").Value = 86666.664
$ws.Range("L67").Value = 86666.664
$ws.Range("N67").Value = -88382.664
$ws.Range("H69").Value = 22197.375
$ws.Range("J69").Value = 50000
$ws.Range("L69").Value = 50000
$ws.Range("N69").Value = -51498
$ws.Range("H72").Value = 22197.375
$ws.Range("J72").Value = 50000
$ws.Range("L72").Value = 150000
$ws.Range("N72").Value = -157488
$ws.Range("H99").Value = 7613.643
$ws.Range("I99").Value = 5519.2
$ws.Range("K99").Value = 5519.2
$ws.Range("M99").Value = -4021.2
$ws.Range("H125").Value = 89618.25
$ws.Range("J125").Value = 89618.25
$ws.Range("L125").Value = 89618.25
$ws.Range("N125").Value = -94538.25
$ws.Range("H126").Value = 7613.643
$ws.Range("I126").Value = 5519.2
$ws.Range("K126").Value = 16557.6
$ws.Range("M126").Value = -14087.6
$ws.Range("H132").Value = 360315.44
$ws.Range("I132").Value = 528138.5600000001
$ws.Range("K132").Value = 1584415.68
$ws.Range("M132").Value = -1581885.68
$ws.Range("H134").Value = 3964
$ws.Range("I134").Value = 2204.389
$ws.Range("K134").Value = 6613.167
$ws.Range("M134").Value = -4078.167
$ws.Range("H136").Value = 8498.637000000001
$ws.Range("I136").Value = 5331.3335
$ws.Range("J136").Value = 12299.4
$ws.Range("K136").Value = 15994.0005
$ws.Range("L136").Value = 36898.2
$ws.Range("M136").Value = -13444.0005
$ws.Range("N136").Value = -41998.2

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 2716.5
$ws.Range("I113").Value = 849
$ws.Range("K113").Value = 2547
$ws.Range("M113").Value = -377
$ws.Range("H131").Value = 4048.2856
$ws.Range("I131").Value = 2232.5
$ws.Range("K131").Value = 6697.5
$ws.Range("M131").Value = -1657.5
$ws.Range("H134").Value = 56013.5
$ws.Range("I134").Value = 38018.332
$ws.Range("K134").Value = 114054.996
$ws.Range("M134").Value = -108984.996

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 4100
$ws.Range("I113").Value = 3542.8572
$ws.Range("J113").Value = 8000
$ws.Range("K113").Value = 3542.8572
$ws.Range("L113").Value = 8000
$ws.Range("M113").Value = -1372.8572
$ws.Range("N113").Value = -12340
$ws.Range("H126").Value = 6539
$ws.Range("I126").Value = 5055.7144
$ws.Range("K126").Value = 15167.1432
$ws.Range("M126").Value = -12697.1432
$ws.Range("H127").Value = 42447.5
$ws.Range("J127").Value = 42447.5
$ws.Range("L127").Value = 42447.5
$ws.Range("N127").Value = -52367.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 6999.2666
$ws.Range("I100").Value = 3925.2856
$ws.Range("J100").Value = 9689
$ws.Range("K100").Value = 3925.2856
$ws.Range("L100").Value = 9689
$ws.Range("M100").Value = -3384.2856
$ws.Range("N100").Value = -10771
$ws.Range("H116").Value = 68000
$ws.Range("J116").Value = 68000
$ws.Range("L116").Value = 68000
$ws.Range("N116").Value = -77178
$ws.Range("H136").Value = 62506840
$ws.Range("I136").Value = 37044390
$ws.Range("K136").Value = 111133170
$ws.Range("M136").Value = -111130620

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H24").Value = 24499.5
$ws.Range("J24").Value = 24499.5
$ws.Range("L24").Value = 24499.5
$ws.Range("N24").Value = -24959.5
$ws.Range("H132").Value = 132738.89
$ws.Range("I132").Value = 172065.27
$ws.Range("J132").Value = 3835.7222
$ws.Range("K132").Value = 516195.8099999999
$ws.Range("L132").Value = 11507.1666
$ws.Range("M132").Value = -513665.8099999999
$ws.Range("N132").Value = -16567.1666
